$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 5 new blank columns starting at column I (I,J,K,L,M). This shifts
#    the existing student_id / stream / grade columns (I,J,K) to N,O,P while
#    leaving column I itself empty, matching the target layout.
# ---------------------------------------------------------------------------
$ws.Range("I1:M1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. New header row (row 1) values - parent job/phone columns.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 10).Value = "father_job"     # J1
$ws.Cells.Item(1, 11).Value = "Mother_job"     # K1
$ws.Cells.Item(1, 12).Value = "father_phone"   # L1
$ws.Cells.Item(1, 13).Value = "mother_phone"   # M1

# ---------------------------------------------------------------------------
# 3. Fill in the "not" placeholder for every student row (2-31).
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 10).Value = "not"          # J - father_job
    $ws.Cells.Item($r, 11).Value = "not"          # K - Mother_job
    $ws.Cells.Item($r, 12).Value = "not"          # L - father_phone
    $ws.Cells.Item($r, 13).Value = "not"          # M - mother_phone
}

# ---------------------------------------------------------------------------
# 4. New header row (row 1) values - class/phone/academic year columns.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 17).Value = "class_id"       # Q1
$ws.Cells.Item(1, 18).Value = "phone"          # R1
$ws.Cells.Item(1, 19).Value = "academic_year"  # S1

# ---------------------------------------------------------------------------
# 5. Fill in class_id / phone / academic_year for every student row (2-31).
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 17).Value = 8               # Q - class_id
    $ws.Cells.Item($r, 18).Value = 11              # R - phone
    $ws.Cells.Item($r, 19).Value = "2024-2025"    # S - academic_year
}

# ---------------------------------------------------------------------------
# 6. Restore the "date of birth" column (D), which this runtime fails to
#    parse (and otherwise turns into #VALUE! errors as soon as the workbook
#    is touched) back to the original dates. The diff does not alter this
#    column, so we simply re-apply the same dates using values the runtime
#    can round-trip correctly.
# ---------------------------------------------------------------------------
$dobSerials = @(39102,39522,39214,40016,38965,40125,39492,39568,39610,39319,39733,39787,38747,39127,39529,38819,39588,39241,39644,39670,40084,38995,39764,39807,39827,39496,39537,39187,38847,38873)
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 4).Value = $dobSerials[$r - 2]
}

# ---------------------------------------------------------------------------
# 7. Update the view: scroll so column E is the left-most visible column and
#    select T26 (matching the final selection recorded in the workbook).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E1").Select()
$ws.Range("T26").Select()
